$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.853.76"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.27%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.636.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.65%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.17%  "

$ws.Range("E6").Value = "  -0.53%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("E8").Value = "  +0.27%  "

$ws.Range("E9").Value = "  +0.89%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.13"
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0783"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.71%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.674.43"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.01%  "

$ws.Range("E13").Value = "  +0.27%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.862.34"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.65%  "

$ws.Range("E15").Value = "  +1.47%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₃0765"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.67%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.27"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.34%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.860.34"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.32%  "

$ws.Range("E19").Value = "  +0.01%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "194.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.28%  "

$ws.Range("E21").Value = "  +1.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.94"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.67%  "

$ws.Range("E23").Value = "  +3.41%  "

$ws.Range("E24").Value = "  +0.07%  "

$ws.Range("E25").Value = "  -3.89%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "138.39"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.06%  "

$ws.Range("E27").Value = "  -4.63%  "

$ws.Range("E28").Value = "  +1.56%  "

$ws.Range("E29").Value = "  +0.86%  "

$ws.Range("E30").Value = "  +0.48%  "

$ws.Range("E31").Value = "  +1.28%  "

$ws.Range("E32").Value = "  +0.57%  "

$ws.Range("E33").Value = "  +1.74%  "

$ws.Range("E34").Value = "  +0.95%  "

$ws.Range("E35").Value = "  +0.79%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.904"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.08%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.58"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.65%  "

$ws.Range("E38").Value = "  +0.18%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.123.81"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.88%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0159"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.78%  "

$ws.Range("E41").Value = "  +0.58%  "

$ws.Range("E42").Value = "  -1.77%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.45"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.26%  "

$ws.Range("E44").Value = "  +1.01%  "

$ws.Range("E45").Value = "  -3.44%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "55.46"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.19%  "

$ws.Range("E47").Value = "  -4.21%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0505"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.35%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.63"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.08%  "

$ws.Range("E50").Value = "  -0.52%  "

$ws.Range("E51").Value = "  -0.06%  "
